# Season record columns: Wins / Losses / Ties
# Adds three new columns (AD, AE, AF) to the existing player-stats sheet:
#   - Header row (row 1) gets labels "Wins", "Losses", "Ties" styled like
#     the other header cells (bold, bordered, centered).
#   - Every data row (2-47) gets the team's season record: 74 wins,
#     87 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 47

# --- Header row ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the preceding header cell (AC1) onto the new
# header cells so they reuse the same header style (bold/border/center).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows ---
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
